$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "COMISIONISTA 4" block (columns I:J) -- this removes the
# header "COMISIONISTA 4" (I1) and its "Porcentaje" column (J1), shifting
# everything to the right of them one-left by two columns.
$ws.Range("I1:J1").EntireColumn.Delete()

# Mirror the selection state left behind by the author's edit (selecting the
# two columns before deleting them / after deletion the same range highlights
# what is now occupied by the former K:L content).
$ws.Range("I1:J1048576").Select()
